# Apply the "Added automation test case changes for ubuntu" edit:
#  - Update region code text from "RR8" to "RR1" in C3, C4, C5
#  - Update B3 label from "Region R1 and Others" to "Derived waiting time"
#  - Align C5's border formatting with C3/C4 (full thin border instead of left/right only)
#  - Move the active selection from E3 to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text changes (set the RR1 values first so new shared strings are appended
# in the same order as the target workbook: RR1 before "Derived waiting time")
$ws.Range("C3").Value = "RR1"
$ws.Range("C4").Value = "RR1"
$ws.Range("C5").Value = "RR1"
$ws.Range("B3").Value = "Derived waiting time"

# Copy C4's cell formatting (full thin border) onto C5 so its border matches
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selected/active cell shown when the workbook is opened
$ws.Range("B9").Select()
